$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix: Property csvEncoding of xava.properties has no effect when
# --- importing CSV files. The accented character now survives the
# --- CSV -> XLSX import, so the previously mangled text is corrected.
$ws.Range("E5").Value = "Playing Fórtnite"

# Normalize the custom date/time number format code to the canonical
# lower-case Excel token casing (m/d/yy h:mm AM/PM) for every cell that
# used the old upper-case custom format.
$ws.Range("A2:A5").NumberFormat = "m/d/yy\ h:mm\ AM/PM"
$ws.Range("D2:D5").NumberFormat = "m/d/yy\ h:mm\ AM/PM"

# Re-save with the current view/selection state.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$null = $ws.Range("E6").Select()

# Header/footer margin, re-expressed from 1.3 cm (was stored with a
# slightly different floating point rounding before).
$ws.PageSetup.HeaderMargin = 36.850393700787386
$ws.PageSetup.FooterMargin = 36.850393700787386
